# Updated symbol list on Mon Jan 16 22:54:12 UTC 2023 with GitHub Actions
# Refresh crypto price/volume snapshot values in Sheet1 (columns D=Price, E=Volume(1h)).
# Values are stored as text (matching the source inlineStr cells), so we force
# NumberFormat to "@" (Text) before assigning each value to stop Excel from
# auto-coercing the numeric-looking strings into Number cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.56%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.91%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.148"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.99%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08099"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "9.81%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.539"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "15.07%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.787"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.62%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.910"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.37%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9333"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.80%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1762"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.36%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07340"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.17%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08871"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "8.96%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03028"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.03%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1000"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001523"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.75%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005782"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.44%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.564"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.07%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.287"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.87%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.28%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1341"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.52%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.158"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-10.46%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.30%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04632"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.04%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001239"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.19%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004522"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.78%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.62%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.46%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01764"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.63%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04603"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.17%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006878"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-4.09%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1374"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.92%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002142"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.83%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.54%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006202"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.20%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7484"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-7.43%"
